$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# remain stored as text (matching the source data, which keeps these
# as inline/shared strings rather than numeric cells).
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D16",
    "D18",
    "D19",
    "D20",
    "D22",
    "D23",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.739.96'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.851.32'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("D4").Value = '1.012'
$ws.Range("E4").Value = '  -2.54%  '
$ws.Range("D5").Value = '319.92'
$ws.Range("E5").Value = '  -1.57%  '
$ws.Range("D6").Value = '1.011'
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("D7").Value = '0.4321'
$ws.Range("E7").Value = '  -2.11%  '
$ws.Range("D8").Value = '0.3764'
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("D9").Value = '0.07387'
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("D10").Value = '0.8852'
$ws.Range("E10").Value = '  +0.12%  '
$ws.Range("D11").Value = '21.69'
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").Value = '1.884.34'
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("D13").Value = '6.759'
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("D14").Value = '5.481'
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("D15").Value = '0.07131'
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("D16").Value = '88.32'
$ws.Range("E16").Value = '  +5.54%  '
$ws.Range("E17").Value = '  -2.24%  '
$ws.Range("D18").Value = '0.000009035'
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").Value = '1.011'
$ws.Range("E19").Value = '  -2.23%  '
$ws.Range("D20").Value = '15.50'
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").Value = '27.739.05'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").Value = '5.277'
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("D23").Value = '11.20'
$ws.Range("E23").Value = '  -1.70%  '
$ws.Range("D24").Value = '2.086.49'
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("D25").Value = '2.032'
$ws.Range("E25").Value = '  +3.35%  '
$ws.Range("D26").Value = '156.01'
$ws.Range("E26").Value = '  -1.44%  '
$ws.Range("D27").Value = '18.59'
$ws.Range("E27").Value = '  -1.47%  '
$ws.Range("D28").Value = '2.125'
$ws.Range("E28").Value = '  +6.82%  '
$ws.Range("D29").Value = '5.425'
$ws.Range("E29").Value = '  +1.90%  '
$ws.Range("D30").Value = '120.99'
$ws.Range("E30").Value = '  +2.95%  '
$ws.Range("D31").Value = '0.08959'
$ws.Range("E31").Value = '  -1.60%  '
$ws.Range("D32").Value = '1.241'
$ws.Range("E32").Value = '  +1.88%  '
$ws.Range("D33").Value = '0.7820'
$ws.Range("E33").Value = '  +0.76%  '
$ws.Range("D34").Value = '4.580'
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("E35").Value = '  -4.99%  '
$ws.Range("D36").Value = '1.148'
$ws.Range("E36").Value = '  -1.47%  '
$ws.Range("D37").Value = '1.012'
$ws.Range("E37").Value = '  -2.35%  '
$ws.Range("D38").Value = '0.05343'
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("E39").Value = '  -1.41%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '7.119'
$ws.Range("E40").Value = '  +3.15%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.863'
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("D42").Value = '0.5203'
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").Value = '0.1687'
$ws.Range("E43").Value = '  -0.52%  '
$ws.Range("D44").Value = '9.003'
$ws.Range("E44").Value = '  +3.24%  '
$ws.Range("D45").Value = '110.89'
$ws.Range("E45").Value = '  +1.23%  '
$ws.Range("D46").Value = '10.78'
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("D47").Value = '1.719'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("D48").Value = '0.4748'
$ws.Range("E48").Value = '  +0.75%  '
$ws.Range("D49").Value = '0.06506'
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("D50").Value = '1.012'
$ws.Range("E50").Value = '  -2.42%  '
$ws.Range("D51").Value = '1.904'
$ws.Range("E51").Value = '  +1.20%  '
